$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 358; this shifts the existing
# rows 358:446 down to 360:448 and keeps their data/formatting intact.
$ws.Rows("358:359").Insert()

# --- New row 358 -----------------------------------------------------
$ws.Range("A358").Value = 4
$ws.Range("B358").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C358").Value = "Los Lagos"
$ws.Range("D358").Value = 44642
$ws.Range("E358").Value = 10
$ws.Range("F358").Value = "Fruta"
$ws.Range("G358").Value = 100102
$ws.Range("H358").Value = "Cítricos"
$ws.Range("I358").Value = 100102003
$ws.Range("J358").Value = "Limón"
$ws.Range("K358").Value = "Sin especificar"
$ws.Range("L358").Value = "1a plateado"
$ws.Range("M358").Value = 1000
$ws.Range("N358").Value = 29000
$ws.Range("O358").Value = 30000
$ws.Range("P358").Value = 29500
$ws.Range("Q358").Value = "`$/malla 18 kilos"
$ws.Range("R358").Value = "Provincia de Melipilla"
$ws.Range("S358").Value = 1639
$ws.Range("T358").Value = 18

# --- New row 359 -----------------------------------------------------
$ws.Range("A359").Value = 4
$ws.Range("B359").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C359").Value = "Los Lagos"
$ws.Range("D359").Value = 44642
$ws.Range("E359").Value = 10
$ws.Range("F359").Value = "Fruta"
$ws.Range("G359").Value = 100102
$ws.Range("H359").Value = "Cítricos"
$ws.Range("I359").Value = 100102003
$ws.Range("J359").Value = "Limón"
$ws.Range("K359").Value = "Sin especificar"
$ws.Range("L359").Value = "2a plateado"
$ws.Range("M359").Value = 500
$ws.Range("N359").Value = 26000
$ws.Range("O359").Value = 26000
$ws.Range("P359").Value = 26000
$ws.Range("Q359").Value = "`$/malla 18 kilos"
$ws.Range("R359").Value = "Provincia de Melipilla"
$ws.Range("S359").Value = 1444
$ws.Range("T359").Value = 18
